# Mexico Liga MX workbook update (28-05-2024)
# - Several adjacent row pairs (and one 3-row group) had their match-data
#   columns (B..AD) swapped between rows while the index column (A) stayed
#   fixed -- the underlying ids/teams/odds were reassigned to the correct
#   row position.
# - A brand new match row (340) was appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param([int]$RowA, [int]$RowB)
    $rA = $ws.Range("B$RowA`:AD$RowA")
    $rB = $ws.Range("B$RowB`:AD$RowB")
    $vA = $rA.Value2
    $vB = $rB.Value2
    $rA.Value = $vB
    $rB.Value = $vA
}

# Simple pairwise swaps (row <-> row), columns B..AD only (column A, the
# running index, is left untouched).
Swap-Rows 36 37
Swap-Rows 94 95
Swap-Rows 128 129
Swap-Rows 148 149
Swap-Rows 175 176
Swap-Rows 200 201
Swap-Rows 222 223
Swap-Rows 264 265
Swap-Rows 272 273
Swap-Rows 276 277
Swap-Rows 298 299
Swap-Rows 318 319
Swap-Rows 322 323

# 3-way rotation: new310 = old311, new311 = old312, new312 = old310
$r310 = $ws.Range("B310:AD310")
$r311 = $ws.Range("B311:AD311")
$r312 = $ws.Range("B312:AD312")
$v310 = $r310.Value2
$v311 = $r311.Value2
$v312 = $r312.Value2
$r310.Value = $v311
$r311.Value = $v312
$r312.Value = $v310

# New row 340 appended at the bottom of the table.
# Copy formatting from the last existing data row (339) for the styled
# columns (A: bold/bordered index style, D: date number format).
$ws.Range("A339").Copy() | Out-Null
$ws.Range("A340").PasteSpecial(-4122) | Out-Null
$ws.Range("D339").Copy() | Out-Null
$ws.Range("D340").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(340, 1).Value = 338
$ws.Cells.Item(340, 2).Value = 8241446
$ws.Cells.Item(340, 3).Value = "Mexico Liga MX"
$ws.Cells.Item(340, 4).Value = 45438.94097222222
$ws.Cells.Item(340, 5).Value = "Club America"
$ws.Cells.Item(340, 6).Value = "Cruz Azul"
$ws.Cells.Item(340, 7).Value = 1
$ws.Cells.Item(340, 8).Value = 0
$ws.Cells.Item(340, 11).Value = "H"
$ws.Cells.Item(340, 12).Value = 1.8
$ws.Cells.Item(340, 13).Value = 3.5
$ws.Cells.Item(340, 14).Value = 4
$ws.Cells.Item(340, 15).Value = 1.85
$ws.Cells.Item(340, 16).Value = 3.4
$ws.Cells.Item(340, 17).Value = 4.2
$ws.Cells.Item(340, 18).Value = -0.5
$ws.Cells.Item(340, 19).Value = 1.925
$ws.Cells.Item(340, 20).Value = 1.925
$ws.Cells.Item(340, 21).Value = 2.25
$ws.Cells.Item(340, 22).Value = 2.025
$ws.Cells.Item(340, 23).Value = 1.825
$ws.Cells.Item(340, 24).Value = 0.8500000000000001
$ws.Cells.Item(340, 25).Value = -1
$ws.Cells.Item(340, 26).Value = -1
$ws.Cells.Item(340, 27).Value = 0.925
$ws.Cells.Item(340, 28).Value = -1
$ws.Cells.Item(340, 29).Value = -1
$ws.Cells.Item(340, 30).Value = 0.825
